$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.265.92"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "1.906.82"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'307.48"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.5243"
$ws.Range("E7").Value = "  +3.37%  "
$ws.Range("D8").Value = "'0.3781"
$ws.Range("E8").Value = "  +3.66%  "
$ws.Range("D9").Value = "'0.07259"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'21.28"
$ws.Range("E10").Value = "  +3.87%  "
$ws.Range("D11").Value = "'0.9008"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "'0.08275"
$ws.Range("E12").Value = "  +10.67%  "
$ws.Range("D13").Value = "1.915.82"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").Value = "'95.41"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "'5.282"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "'0.000008606"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "'14.49"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "'1.0000"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "27.289.97"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "2.158.63"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'10.66"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").Value = "'6.461"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = "  +10.58%  "
$ws.Range("D26").Value = "'146.00"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "'114.84"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("E30").Value = "  +6.42%  "
$ws.Range("D31").Value = "'4.817"
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("D32").Value = "'0.09208"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").Value = "'0.8060"
$ws.Range("E33").Value = "  +8.17%  "
$ws.Range("D34").Value = "'0.05086"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("D35").Value = "'1.242"
$ws.Range("E35").Value = "  +8.34%  "
$ws.Range("D36").Value = "'2.954"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'3.339"
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("D38").Value = "'2.574"
$ws.Range("E38").Value = "  +3.18%  "
$ws.Range("D39").Value = "'0.5733"
$ws.Range("E39").Value = "  +3.82%  "
$ws.Range("D40").Value = "'0.01979"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'1.076"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'9.073"
$ws.Range("E42").Value = "  +6.33%  "
$ws.Range("D43").Value = "'6.637"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").Value = "'118.77"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("D45").Value = "'0.1518"
$ws.Range("E45").Value = "  +2.64%  "
$ws.Range("D46").Value = "'0.4846"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'10.16"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "'1.613"
$ws.Range("E49").Value = "  +4.28%  "
$ws.Range("D50").Value = "'37.61"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").Value = "'63.68"
$ws.Range("E51").Value = "  +1.47%  "
